$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.365.67'
$ws.Range("E2").Value = '  +2.56%  '

$ws.Range("D3").Value = '1.838.98'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.44%  '

$ws.Range("E6").Value = '  +1.36%  '

$ws.Range("E7").Value = '  +0.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.11'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +12.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.307'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0698'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.32%  '

$ws.Range("E11").Value = '  +2.96%  '

$ws.Range("D12").Value = '2.107.49'
$ws.Range("E12").Value = '  +1.90%  '

$ws.Range("D13").Value = '1.836.46'
$ws.Range("E13").Value = '  +1.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.668'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.61%  '

$ws.Range("D17").Value = '35.345.23'
$ws.Range("E17").Value = '  +2.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("D19").Value = '0.0₃0797'
$ws.Range("E19").Value = '  +4.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '244.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.35%  '

$ws.Range("E22").Value = '  +14.22%  '

$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '169.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.16%  '

$ws.Range("E28").Value = '  +0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +25.17%  '

$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("D31").Value = '3.260.18'
$ws.Range("E31").Value = '  +34.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0547'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.34%  '

$ws.Range("E34").Value = '  +5.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.73%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '93.66'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.682'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.82%  '

$ws.Range("D38").Value = '1.341.35'
$ws.Range("E38").Value = '  +2.30%  '

$ws.Range("E39").Value = '  +2.80%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0194'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.40%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.25%  '

$ws.Range("E43").Value = '  +5.77%  '

$ws.Range("E44").Value = '  +3.83%  '

$ws.Range("E45").Value = '  +0.84%  '

$ws.Range("E46").Value = '  +0.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.38%  '

$ws.Range("E48").Value = '  +0.98%  '

$ws.Range("D49").Value = '2.008.11'
$ws.Range("E49").Value = '  +1.98%  '

$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.00%  '
